# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" worksheet (a copy of the existing quarterly
#    report sheet, so it inherits identical sheetPr/pageMargins/styles)
#    positioned right after the "总计" summary sheet.
# 2. Populate the new sheet with the 2022-Q3 fund-holding table.
# 3. Update the "总计" summary sheet: insert a new data row for 2022-Q3
#    at the top of the table and renumber the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the current second sheet (2022-Q2) so the new sheet
# inherits the same sheetPr / styles / pageMargins, then move it into
# place and rename it.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(2)
$templateSheet.Copy($templateSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Clear the copied sample data (keep header row styling/layout intact)
$newSheet.Range("A2:H4").ClearContents()

# ---------------------------------------------------------------------
# Step 2: write the 2022-Q3 fund table.
# Columns D/E/F/G and the fund-code column B hold numeric-looking text
# in this workbook (e.g. "25.18", "012526") so a leading apostrophe is
# used to force them to be stored as text instead of numbers. Column H
# (仓位排名) is a genuine number.
# ---------------------------------------------------------------------
$rows = @(
    @("012526", "广发盛锦混合型证券投资基金A", "25.18", "90.19", "4.94", "1.2439", 4),
    @("012527", "广发盛锦混合型证券投资基金C", "1.16", "90.19", "4.94", "0.0573", 4),
    @("005075", "富国研究量化精选混合", "2.60", "90.00", "1.62", "0.0421", 9),
    @("009719", "招商增浩一年定期开放混合C", "1.34", "23.60", "0.66", "0.0088", 9),
    @("003456", "信澳新目标灵活配置混合", "0.39", "94.17", "1.27", "0.0050", 10),
    @("009718", "招商增浩一年定期开放混合A", "0.70", "23.60", "0.66", "0.0046", 9),
    @("002952", "建信多因子量化股票", "0.09", "91.26", "3.01", "0.0027", 10)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $newSheet.Range("A$r").Value = $i
    $newSheet.Range("B$r").Value = "'" + $data[0]
    $newSheet.Range("C$r").Value = $data[1]
    $newSheet.Range("D$r").Value = "'" + $data[2]
    $newSheet.Range("E$r").Value = "'" + $data[3]
    $newSheet.Range("F$r").Value = "'" + $data[4]
    $newSheet.Range("G$r").Value = "'" + $data[5]
    $newSheet.Range("H$r").Value = $data[6]
}

# The A column (row index) carries the same bold/border style throughout
# the table; make sure the newly-added rows (5-8) pick it up too by
# copying the format from the existing A2 cell.
$newSheet.Range("A2").Copy()
$newSheet.Range("A5:A8").PasteSpecial(-4122)
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $newSheet.Range("A$r").Value = $i
}

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet with the new 2022-Q3 row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Make room for the new row 6 (style for the new A6 index cell is
# copied from the current last row, A5, before it gets overwritten).
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 1.36

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.48

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 4
$totalSheet.Range("D4").Value = 0.47

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 3
$totalSheet.Range("D5").Value = 0.19

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q3"
$totalSheet.Range("C6").Value = 1
$totalSheet.Range("D6").Value = 0.02

# ---------------------------------------------------------------------
# Keep the trailing "2021-Q3" sheet as the active tab (matches the
# original workbook's tab-selection state).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
